$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 26.81310566666667
$ws.Range("H2").Value2 = 80.439317
$ws.Range("I2").Value2 = 0.004518206005002021
$ws.Range("J2").Value2 = 0.004518206005002021
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 7.107333666666666
$ws.Range("N2").Value2 = 21.322001
$ws.Range("O2").Value2 = 0.7373665550576455
$ws.Range("P2").Value2 = 0.7373665550576454
$ws.Range("Q2").Value2 = 190.5696886125908
$ws.Range("R2").Value2 = 1715.127197513317
$ws.Range("S2").Value2 = 0.003331573996949107
$ws.Range("T2").Value2 = 0.003331573996949107

# Row 3
$ws.Range("G3").Value2 = 26.81310566666667
$ws.Range("H3").Value2 = 80.439317
$ws.Range("I3").Value2 = 0.004518206005002021
$ws.Range("J3").Value2 = 0.004518206005002021
$ws.Range("O3").Value2 = 0.1688878844614928
$ws.Range("P3").Value2 = 0.1688878844614928
$ws.Range("Q3").Value2 = 43.64845588874011
$ws.Range("R3").Value2 = 392.836102998661
$ws.Range("S3").Value2 = 0.0007630702537460044
$ws.Range("T3").Value2 = 0.0007630702537460042

# Row 4
$ws.Range("G4").Value2 = 26.81310566666667
$ws.Range("H4").Value2 = 80.439317
$ws.Range("I4").Value2 = 0.004518206005002021
$ws.Range("J4").Value2 = 0.004518206005002021
$ws.Range("M4").Value2 = 0.8135026666666666
$ws.Range("N4").Value2 = 2.440508
$ws.Range("O4").Value2 = 0.08439869112428164
$ws.Range("P4").Value2 = 0.08439869112428162
$ws.Range("Q4").Value2 = 21.81253296144844
$ws.Range("R4").Value2 = 196.312796653036
$ws.Range("S4").Value2 = 0.00038133067305204
$ws.Range("T4").Value2 = 0.00038133067305204

# Row 5
$ws.Range("G5").Value2 = 26.81310566666667
$ws.Range("H5").Value2 = 80.439317
$ws.Range("I5").Value2 = 0.004518206005002021
$ws.Range("J5").Value2 = 0.004518206005002021
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.09009266666666667
$ws.Range("N5").Value2 = 0.270278
$ws.Range("O5").Value2 = 0.009346869356580103
$ws.Range("P5").Value2 = 0.009346869356580103
$ws.Range("Q5").Value2 = 2.415664191125111
$ws.Range("R5").Value2 = 21.740977720126
$ws.Range("S5").Value2 = 0.0000422310812548696
$ws.Range("T5").Value2 = 0.0000422310812548696

# Row 6
$ws.Range("G6").Value2 = 5771.873535333333
$ws.Range("I6").Value2 = 0.9726032482643521
$ws.Range("J6").Value2 = 0.9726032482643523
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 7.107333666666666
$ws.Range("N6").Value2 = 21.322001
$ws.Range("O6").Value2 = 0.7373665550576455
$ws.Range("P6").Value2 = 0.7373665550576454
$ws.Range("Q6").Value2 = 41022.63109741695
$ws.Range("R6").Value2 = 369203.6798767526
$ws.Range("S6").Value2 = 0.7171651066105612
$ws.Range("T6").Value2 = 0.7171651066105612

# Row 7
$ws.Range("G7").Value2 = 5771.873535333333
$ws.Range("I7").Value2 = 0.9726032482643521
$ws.Range("J7").Value2 = 0.9726032482643523
$ws.Range("O7").Value2 = 0.1688878844614928
$ws.Range("P7").Value2 = 0.1688878844614928
$ws.Range("Q7").Value2 = 9395.904022993511
$ws.Range("R7").Value2 = 84563.13620694159
$ws.Range("S7").Value2 = 0.1642609050197425
$ws.Range("T7").Value2 = 0.1642609050197425

# Row 8
$ws.Range("G8").Value2 = 5771.873535333333
$ws.Range("I8").Value2 = 0.9726032482643521
$ws.Range("J8").Value2 = 0.9726032482643523
$ws.Range("M8").Value2 = 0.8135026666666666
$ws.Range("N8").Value2 = 2.440508
$ws.Range("O8").Value2 = 0.08439869112428164
$ws.Range("P8").Value2 = 0.08439869112428162
$ws.Range("Q8").Value2 = 4695.434512656427
$ws.Range("R8").Value2 = 42258.91061390785
$ws.Range("S8").Value2 = 0.08208644113673606
$ws.Range("T8").Value2 = 0.08208644113673606

# Row 9
$ws.Range("G9").Value2 = 5771.873535333333
$ws.Range("I9").Value2 = 0.9726032482643521
$ws.Range("J9").Value2 = 0.9726032482643523
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 0.6666666666666666
$ws.Range("M9").Value2 = 0.09009266666666667
$ws.Range("N9").Value2 = 0.270278
$ws.Range("O9").Value2 = 0.009346869356580103
$ws.Range("P9").Value2 = 0.009346869356580103
$ws.Range("Q9").Value2 = 520.0034784609409
$ws.Range("R9").Value2 = 4680.031306148469
$ws.Range("S9").Value2 = 0.009090795497312343
$ws.Range("T9").Value2 = 0.009090795497312345

# Row 10
$ws.Range("G10").Value2 = 132.4457753333333
$ws.Range("H10").Value2 = 397.337326
$ws.Range("I10").Value2 = 0.02231808970163987
$ws.Range("J10").Value2 = 0.02231808970163988
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 7.107333666666666
$ws.Range("N10").Value2 = 21.322001
$ws.Range("O10").Value2 = 0.7373665550576455
$ws.Range("P10").Value2 = 0.7373665550576454
$ws.Range("Q10").Value2 = 941.3363180343695
$ws.Range("R10").Value2 = 8472.026862309327
$ws.Range("S10").Value2 = 0.01645661291876571
$ws.Range("T10").Value2 = 0.01645661291876571

# Row 11
$ws.Range("G11").Value2 = 132.4457753333333
$ws.Range("H11").Value2 = 397.337326
$ws.Range("I11").Value2 = 0.02231808970163987
$ws.Range("J11").Value2 = 0.02231808970163988
$ws.Range("O11").Value2 = 0.1688878844614928
$ws.Range("P11").Value2 = 0.1688878844614928
$ws.Range("Q11").Value2 = 215.6055197094842
$ws.Range("R11").Value2 = 1940.449677385358
$ws.Range("S11").Value2 = 0.003769254954931788
$ws.Range("T11").Value2 = 0.003769254954931787

# Row 12
$ws.Range("G12").Value2 = 132.4457753333333
$ws.Range("H12").Value2 = 397.337326
$ws.Range("I12").Value2 = 0.02231808970163987
$ws.Range("J12").Value2 = 0.02231808970163988
$ws.Range("M12").Value2 = 0.8135026666666666
$ws.Range("N12").Value2 = 2.440508
$ws.Range("O12").Value2 = 0.08439869112428164
$ws.Range("P12").Value2 = 0.08439869112428162
$ws.Range("Q12").Value2 = 107.7449914224009
$ws.Range("R12").Value2 = 969.7049228016079
$ws.Range("S12").Value2 = 0.001883617559212715
$ws.Range("T12").Value2 = 0.001883617559212715

# Row 13
$ws.Range("G13").Value2 = 132.4457753333333
$ws.Range("H13").Value2 = 397.337326
$ws.Range("I13").Value2 = 0.02231808970163987
$ws.Range("J13").Value2 = 0.02231808970163988
$ws.Range("K13").Value2 = 2
$ws.Range("L13").Value2 = 0.6666666666666666
$ws.Range("M13").Value2 = 0.09009266666666667
$ws.Range("N13").Value2 = 0.270278
$ws.Range("O13").Value2 = 0.009346869356580103
$ws.Range("P13").Value2 = 0.009346869356580103
$ws.Range("Q13").Value2 = 11.93239308851422
$ws.Range("R13").Value2 = 107.391537796628
$ws.Range("S13").Value2 = 0.0002086042687296637
$ws.Range("T13").Value2 = 0.0002086042687296637

# Row 14
$ws.Range("G14").Value2 = 3.326003
$ws.Range("H14").Value2 = 9.978009
$ws.Range("I14").Value2 = 0.0005604560290058679
$ws.Range("J14").Value2 = 0.000560456029005868
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 7.107333666666666
$ws.Range("N14").Value2 = 21.322001
$ws.Range("O14").Value2 = 0.7373665550576455
$ws.Range("P14").Value2 = 0.7373665550576454
$ws.Range("Q14").Value2 = 23.63901309733433
$ws.Range("R14").Value2 = 212.751117876009
$ws.Range("S14").Value2 = 0.0004132615313693446
$ws.Range("T14").Value2 = 0.0004132615313693446

# Row 15
$ws.Range("G15").Value2 = 3.326003
$ws.Range("H15").Value2 = 9.978009
$ws.Range("I15").Value2 = 0.0005604560290058679
$ws.Range("J15").Value2 = 0.000560456029005868
$ws.Range("O15").Value2 = 0.1688878844614928
$ws.Range("P15").Value2 = 0.1688878844614928
$ws.Range("Q15").Value2 = 5.414326002966333
$ws.Range("R15").Value2 = 48.728934026697
$ws.Range("S15").Value2 = 0.00009465423307249009
$ws.Range("T15").Value2 = 0.00009465423307249008

# Row 16
$ws.Range("G16").Value2 = 3.326003
$ws.Range("H16").Value2 = 9.978009
$ws.Range("I16").Value2 = 0.0005604560290058679
$ws.Range("J16").Value2 = 0.000560456029005868
$ws.Range("M16").Value2 = 0.8135026666666666
$ws.Range("N16").Value2 = 2.440508
$ws.Range("O16").Value2 = 0.08439869112428164
$ws.Range("P16").Value2 = 0.08439869112428162
$ws.Range("Q16").Value2 = 2.705712309841333
$ws.Range("R16").Value2 = 24.351410788572
$ws.Range("S16").Value2 = 0.00004730175528080768
$ws.Range("T16").Value2 = 0.00004730175528080768

# Row 17
$ws.Range("G17").Value2 = 3.326003
$ws.Range("H17").Value2 = 9.978009
$ws.Range("I17").Value2 = 0.0005604560290058679
$ws.Range("J17").Value2 = 0.000560456029005868
$ws.Range("K17").Value2 = 2
$ws.Range("L17").Value2 = 0.6666666666666666
$ws.Range("M17").Value2 = 0.09009266666666667
$ws.Range("N17").Value2 = 0.270278
$ws.Range("O17").Value2 = 0.009346869356580103
$ws.Range("P17").Value2 = 0.009346869356580103
$ws.Range("Q17").Value2 = 0.2996484796113333
$ws.Range("R17").Value2 = 2.696836316502
$ws.Range("S17").Value2 = 0.000005238509283225516
$ws.Range("T17").Value2 = 0.000005238509283225517
